$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet2: mirror Sheet1's A:B data (preserves types/styles/shared strings),
#     then add C:D columns with the new "same"/"expected" test columns ---
$ws1.Range("A1:B3").Copy($ws2.Range("A1"))
$ws1.Range("A5:B5").Copy($ws2.Range("A5"))

# Header row for the new columns
$ws2.Range("C1").Value = "same"
$ws2.Range("D1").Value = "expected"

# Row 2
$ws2.Range("C2").Value = "b"
$ws2.Range("D2").Value = "b"

# Row 3
$ws2.Range("C3").Value = "b"
$ws2.Range("D3").Value = $false

# Row 5
$ws2.Range("C5").Value = "b"
$ws2.Range("D5").Value = "howdy"

# Set the desired final selection on Sheet2, then re-select the Sheet1
# range so Sheet1 remains the active (tabSelected) sheet, matching the diff.
$ws2.Range("A1:D5").Select() | Out-Null
$ws1.Range("C1:D5").Select() | Out-Null
